# Add "PVTStL" (Policy Vehicle Types Subject to LCFS) to the acronym key
# on the "Key to Variables" sheet, inserted as a new row 250 (pushing the
# existing rows 250-271 down to 251-272).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")
$ws.Activate()

# Insert a brand-new row at position 250; Excel shifts rows 250.. down
# and copies formatting from the row above, matching the style pattern
# (A/B/C = fill style 2, F = fill style 4) used by neighboring "trans" rows.
$ws.Rows.Item(250).Insert()

$ws.Cells.Item(250, 1).Value = "trans"
$ws.Cells.Item(250, 2).Value = "PVTStL"
$ws.Cells.Item(250, 3).Value = "Policy Vehicle Types Subject to LCFS"
$ws.Cells.Item(250, 6).Value = "low"

# Mirror the author's final selection / scroll position on the sheet.
$ws.Range("A250:F250").Select()
$excel.ActiveWindow.ScrollRow = 239
